$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 216.33333
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 274.5
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 274.5
$ws.Range("M9").Value = 69
$ws.Range("N9").Value = -612.5
$ws.Range("H15").Value = 1106.3422
$ws.Range("I15").Value = 1106.3422
$ws.Range("K15").Value = 3319.0266
$ws.Range("M15").Value = -3150.0266
$ws.Range("H18").Value = 888.7
$ws.Range("I18").Value = 896
$ws.Range("K18").Value = 896
$ws.Range("M18").Value = -612
$ws.Range("H116").Value = 4610.6313
$ws.Range("J116").Value = 5685.25
$ws.Range("L116").Value = 5685.25
$ws.Range("N116").Value = -12569.25
$ws.Range("H131").Value = 818.8
$ws.Range("I131").Value = 818.8
$ws.Range("K131").Value = 2456.4
$ws.Range("M131").Value = 2583.6
$ws.Range("H132").Value = 1990.2987
$ws.Range("J132").Value = 1662
$ws.Range("L132").Value = 4986
$ws.Range("N132").Value = -10046
$ws.Range("H138").Value = 2830.2778
$ws.Range("I138").Value = 1266.1666
$ws.Range("J138").Value = 5020.033
$ws.Range("K138").Value = 3798.4998
$ws.Range("L138").Value = 15060.099
$ws.Range("M138").Value = 1341.5002
$ws.Range("N138").Value = -25340.099
$ws.Range("H141").Value = 688.3684
$ws.Range("I141").Value = 687.94446
$ws.Range("K141").Value = 2063.83338
$ws.Range("M141").Value = 3116.16662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7792.567
$ws.Range("I32").Value = 1036.625
$ws.Range("K32").Value = 1036.625
$ws.Range("M32").Value = -749.625
$ws.Range("H74").Value = 1790.875
$ws.Range("I74").Value = 1507.9259
$ws.Range("J74").Value = 3318.8
$ws.Range("K74").Value = 1507.9259
$ws.Range("L74").Value = 3318.8
$ws.Range("M74").Value = -633.9259
$ws.Range("N74").Value = -5066.8
$ws.Range("H77").Value = 1790.875
$ws.Range("I77").Value = 1507.9259
$ws.Range("J77").Value = 3318.8
$ws.Range("K77").Value = 7539.6295
$ws.Range("L77").Value = 16594
$ws.Range("M77").Value = -3171.6295
$ws.Range("N77").Value = -25330
$ws.Range("H97").Value = 430.84848
$ws.Range("I97").Value = 452.55554
$ws.Range("K97").Value = 452.55554
$ws.Range("M97").Value = 43.44445999999999
$ws.Range("H132").Value = 6285
$ws.Range("I132").Value = 5999
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 17997
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -15467
$ws.Range("N132").Value = -25059.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8083411
$ws.Range("I94").Value = 14706984
$ws.Range("K94").Value = 14706984
$ws.Range("M94").Value = -14706533
$ws.Range("H107").Value = 2743.4443
$ws.Range("I107").Value = 2102
$ws.Range("J107").Value = 4026.3333
$ws.Range("K107").Value = 2102
$ws.Range("L107").Value = 4026.3333
$ws.Range("M107").Value = -182
$ws.Range("N107").Value = -7866.3333
$ws.Range("H134").Value = 2324.8635
$ws.Range("I134").Value = 2312.795
$ws.Range("K134").Value = 6938.385
$ws.Range("M134").Value = -4403.385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 550.3333
$ws.Range("I22").Value = 515.625
$ws.Range("J22").Value = 590
$ws.Range("K22").Value = 515.625
$ws.Range("L22").Value = 590
$ws.Range("M22").Value = -165.625
$ws.Range("N22").Value = -1290
$ws.Range("H122").Value = 233419.16
$ws.Range("J122").Value = 1497.5
$ws.Range("L122").Value = 4492.5
$ws.Range("N122").Value = -9392.5
$ws.Range("H134").Value = 1504.6428
$ws.Range("I134").Value = 1516.62
$ws.Range("J134").Value = 1404.8334
$ws.Range("K134").Value = 4549.86
$ws.Range("L134").Value = 4214.5002
$ws.Range("M134").Value = -2014.86
$ws.Range("N134").Value = -9284.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 634.1818
$ws.Range("I103").Value = 964.2
$ws.Range("J103").Value = 359.16666
$ws.Range("K103").Value = 2892.6
$ws.Range("L103").Value = 1077.49998
$ws.Range("M103").Value = -2013.6
$ws.Range("N103").Value = -2835.49998
$ws.Range("H107").Value = 2100.75
$ws.Range("I107").Value = 2133
$ws.Range("J107").Value = 2004
$ws.Range("K107").Value = 6399
$ws.Range("L107").Value = 6012
$ws.Range("M107").Value = -4479
$ws.Range("N107").Value = -9852
$ws.Range("H113").Value = 2086.3125
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 2162.0667
$ws.Range("K113").Value = 2850
$ws.Range("L113").Value = 6486.2001
$ws.Range("M113").Value = -680
$ws.Range("N113").Value = -10826.2001
$ws.Range("H117").Value = 6503.6665
$ws.Range("I117").Value = 4008
$ws.Range("J117").Value = 8999.333000000001
$ws.Range("K117").Value = 12024
$ws.Range("L117").Value = 26997.999
$ws.Range("M117").Value = -8582
$ws.Range("N117").Value = -33881.999
$ws.Range("H120").Value = 7495
$ws.Range("I120").Value = 7495
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 22485
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -17647
$ws.Range("N120").ClearContents()
$ws.Range("H136").Value = 6286.0557
$ws.Range("I136").Value = 5543
$ws.Range("K136").Value = 16629
$ws.Range("M136").Value = -11529
$ws.Range("H137").Value = 3764.9583
$ws.Range("I137").Value = 2267.4211
$ws.Range("J137").Value = 9455.6
$ws.Range("K137").Value = 6802.263300000001
$ws.Range("L137").Value = 28366.8
$ws.Range("M137").Value = -1702.263300000001
$ws.Range("N137").Value = -38566.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4428879
$ws.Range("I3").Value = 7500176
$ws.Range("J3").Value = 333816.66
$ws.Range("K3").Value = 7500176
$ws.Range("L3").Value = 333816.66
$ws.Range("M3").Value = -7500060
$ws.Range("N3").Value = -334048.66
$ws.Range("H25").Value = 3366.1667
$ws.Range("I25").Value = 1349.5
$ws.Range("J25").Value = 7399.5
$ws.Range("K25").Value = 1349.5
$ws.Range("L25").Value = 7399.5
$ws.Range("M25").Value = -820.5
$ws.Range("N25").Value = -8457.5
$ws.Range("H132").Value = 5471.1333
$ws.Range("I132").Value = 5471.1333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16413.3999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13883.3999
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 166075.38
$ws.Range("J135").Value = 88785.30499999999
$ws.Range("L135").Value = 88785.30499999999
$ws.Range("N135").Value = -98925.30499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 8016
$ws.Range("I30").Value = 8016
$ws.Range("K30").Value = 8016
$ws.Range("M30").Value = -7908
$ws.Range("H55").Value = 1399.7894
$ws.Range("I55").Value = 1070.3529
$ws.Range("J55").Value = 4200
$ws.Range("K55").Value = 1070.3529
$ws.Range("L55").Value = 4200
$ws.Range("M55").Value = -897.3529000000001
$ws.Range("N55").Value = -4546
$ws.Range("H122").Value = 4398.811
$ws.Range("I122").Value = 3569.7812
$ws.Range("J122").Value = 9704.6
$ws.Range("K122").Value = 10709.3436
$ws.Range("L122").Value = 29113.8
$ws.Range("M122").Value = -8259.3436
$ws.Range("N122").Value = -34013.8
$ws.Range("H132").Value = 16310.2
$ws.Range("I132").Value = 14112.75
$ws.Range("K132").Value = 42338.25
$ws.Range("M132").Value = -39808.25
$ws.Range("H136").Value = 1577.9111
$ws.Range("I136").Value = 1607.1163
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 4821.3489
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -2271.3489
$ws.Range("N136").Value = -7950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3134.5
$ws.Range("J96").Value = 3780
$ws.Range("L96").Value = 3780
$ws.Range("N96").Value = -6526
$ws.Range("H100").Value = 1837
$ws.Range("I100").Value = 1996.25
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 3992.5
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -3451.5
$ws.Range("N100").Value = -3482
$ws.Range("H132").Value = 2533.6606
$ws.Range("I132").Value = 2453.7754
$ws.Range("J132").Value = 3092.8572
$ws.Range("K132").Value = 7361.3262
$ws.Range("L132").Value = 9278.571599999999
$ws.Range("M132").Value = -4831.3262
$ws.Range("N132").Value = -14338.5716
$ws.Range("H136").Value = 3893.2903
$ws.Range("I136").Value = 4035.875
$ws.Range("J136").Value = 3404.4285
$ws.Range("K136").Value = 12107.625
$ws.Range("L136").Value = 10213.2855
$ws.Range("M136").Value = -9557.625
$ws.Range("N136").Value = -15313.2855

Write-Host "Applied all updates."